$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, center-top alignment) from H1 onto the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New numeric data for columns I (I0) and J (IF), rows 2-84
$ijData = @{
    2 = @(7, 7)
    3 = @(9, 9)
    4 = @(7, 8)
    5 = @(8, 8)
    6 = @(9, 9)
    7 = @(7, 8)
    8 = @(7, 8)
    9 = @(8, 8)
    10 = @(7, 8)
    11 = @(6, 6)
    12 = @(9, 9)
    13 = @(7, 7)
    14 = @(7, 8)
    15 = @(7, 8)
    16 = @(6, 7)
    17 = @(8, 9)
    18 = @(7, 8)
    19 = @(8, 8)
    20 = @(6, 7)
    21 = @(6, 7)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(9, 9)
    25 = @(9, 9)
    26 = @(7, 7)
    27 = @(8, 8)
    28 = @(7, 8)
    29 = @(7, 7)
    30 = @(9, 9)
    31 = @(8, 9)
    32 = @(6, 7)
    33 = @(6, 6)
    34 = @(6, 6)
    35 = @(10, 10)
    36 = @(7, 7)
    37 = @(9, 10)
    38 = @(8, 8)
    39 = @(6, 7)
    40 = @(7, 7)
    41 = @(8, 8)
    42 = @(6, 7)
    43 = @(2, 3)
    44 = @(7, 7)
    45 = @(5, 5)
    46 = @(8, 8)
    47 = @(7, 7)
    48 = @(8, 8)
    49 = @(9, 9)
    50 = @(7, 7)
    51 = @(7, 7)
    52 = @(7, 7)
    53 = @(6, 6)
    54 = @(7, 7)
    55 = @(6, 7)
    56 = @(8, 8)
    57 = @(6, 6)
    58 = @(6, 7)
    59 = @(7, 7)
    60 = @(8, 8)
    61 = @(8, 9)
    62 = @(7, 7)
    63 = @(8, 8)
    64 = @(7, 7)
    65 = @(9, 9)
    66 = @(6, 6)
    67 = @(7, 7)
    68 = @(7, 7)
    69 = @(7, 8)
    70 = @(7, 8)
    71 = @(9, 9)
    72 = @(6, 6)
    73 = @(8, 9)
    74 = @(7, 7)
    75 = @(6, 6)
    76 = @(7, 8)
    77 = @(7, 7)
    78 = @(6, 6)
    79 = @(9, 9)
    80 = @(8, 8)
    81 = @(6, 6)
    82 = @(5, 5)
    83 = @(4, 4)
    84 = @(2, 2)
}

foreach ($row in $ijData.Keys) {
    $vals = $ijData[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
